$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "36.280.41"
$rng.Style = "Normal"

$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = "  -3.17%  "
$rng.Style = "Normal"

$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "1.976.63"
$rng.Style = "Normal"

$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "  -3.83%  "
$rng.Style = "Normal"

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = "  +0.05%  "
$rng.Style = "Normal"

$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "244.82"
$rng.Style = "Normal"

$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = "  -3.32%  "
$rng.Style = "Normal"

$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "0.620"
$rng.Style = "Normal"

$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$rng.Value = "  -5.05%  "
$rng.Style = "Normal"

$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "59.21"
$rng.Style = "Normal"

$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$rng.Value = "  -10.11%  "
$rng.Style = "Normal"

$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = "  +0.06%  "
$rng.Style = "Normal"

$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.374"
$rng.Style = "Normal"

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "  -2.47%  "
$rng.Style = "Normal"

$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "57.15"
$rng.Style = "Normal"

$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "  -4.03%  "
$rng.Style = "Normal"

$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.0828"
$rng.Style = "Normal"

$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = "  +7.59%  "
$rng.Style = "Normal"

$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = "  -0.93%  "
$rng.Style = "Normal"

$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "23.54"
$rng.Style = "Normal"

$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = "  +5.05%  "
$rng.Style = "Normal"

$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "0.861"
$rng.Style = "Normal"

$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = "  -7.01%  "
$rng.Style = "Normal"

$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "2.269.68"
$rng.Style = "Normal"

$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = "  -3.67%  "
$rng.Style = "Normal"

$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "13.91"
$rng.Style = "Normal"

$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = "  -6.61%  "
$rng.Style = "Normal"

$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "5.45"
$rng.Style = "Normal"

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = "  -2.36%  "
$rng.Style = "Normal"

$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "1.980.56"
$rng.Style = "Normal"

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = "  -3.57%  "
$rng.Style = "Normal"

$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "36.217.69"
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "  -3.06%  "
$rng.Style = "Normal"

$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "69.97"
$rng.Style = "Normal"

$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = "  -5.00%  "
$rng.Style = "Normal"

$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "0.0₃0866"
$rng.Style = "Normal"

$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "  -1.07%  "
$rng.Style = "Normal"

$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "5.30"
$rng.Style = "Normal"

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = "  -3.84%  "
$rng.Style = "Normal"

$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "234.28"
$rng.Style = "Normal"

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "  -2.35%  "
$rng.Style = "Normal"

$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = "  -0.11%  "
$rng.Style = "Normal"

$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = "  -2.86%  "
$rng.Style = "Normal"

$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = "  -4.10%  "
$rng.Style = "Normal"

$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "10.00"
$rng.Style = "Normal"

$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = "  +0.07%  "
$rng.Style = "Normal"

$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "161.97"
$rng.Style = "Normal"

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "  -0.12%  "
$rng.Style = "Normal"

$rng = $ws.Range("B29")
$rng.NumberFormat = "@"
$rng.Value = "Kaspa"
$rng.Style = "Normal"

$rng = $ws.Range("C29")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$rng.Style = "Normal"

$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "0.132"
$rng.Style = "Normal"

$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "  +7.53%  "
$rng.Style = "Normal"

$rng = $ws.Range("B30")
$rng.NumberFormat = "@"
$rng.Value = "EthereumClassic"
$rng.Style = "Normal"

$rng = $ws.Range("C30")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$rng.Style = "Normal"

$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "19.75"
$rng.Style = "Normal"

$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = "  -1.33%  "
$rng.Style = "Normal"

$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = "  -1.86%  "
$rng.Style = "Normal"

$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "1.18"
$rng.Style = "Normal"

$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = "  -3.32%  "
$rng.Style = "Normal"

$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "4.89"
$rng.Style = "Normal"

$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = "  -7.23%  "
$rng.Style = "Normal"

$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = "  -0.33%  "
$rng.Style = "Normal"

$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = "  -6.16%  "
$rng.Style = "Normal"

$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = "  +2.73%  "
$rng.Style = "Normal"

$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = "  +0.09%  "
$rng.Style = "Normal"

$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = "  -8.63%  "
$rng.Style = "Normal"

$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = "  -2.25%  "
$rng.Style = "Normal"

$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "3.02"
$rng.Style = "Normal"

$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = "  -0.82%  "
$rng.Style = "Normal"

$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = "  -0.47%  "
$rng.Style = "Normal"

$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "0.0965"
$rng.Style = "Normal"

$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "  -6.69%  "
$rng.Style = "Normal"

$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = "  -4.46%  "
$rng.Style = "Normal"

$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = "  -2.37%  "
$rng.Style = "Normal"

$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = "  -5.44%  "
$rng.Style = "Normal"

$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "16.19"
$rng.Style = "Normal"

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = "  -8.97%  "
$rng.Style = "Normal"

$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "92.25"
$rng.Style = "Normal"

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "  -5.03%  "
$rng.Style = "Normal"

$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "1.367.94"
$rng.Style = "Normal"

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = "  -3.64%  "
$rng.Style = "Normal"

$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "7.48"
$rng.Style = "Normal"

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = "  -5.74%  "
$rng.Style = "Normal"

$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = "  -3.62%  "
$rng.Style = "Normal"

$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "45.41"
$rng.Style = "Normal"

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = "  -3.11%  "
$rng.Style = "Normal"

Write-Host "Applied all cryptos updates"
